$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final model-name ordering for rows 2..26 (row 1 is the header, unchanged)
$names = @(
    "model_11_5_0",
    "model_11_5_22",
    "model_11_5_21",
    "model_11_5_20",
    "model_11_5_19",
    "model_11_5_18",
    "model_11_5_17",
    "model_11_5_16",
    "model_11_5_15",
    "model_11_5_14",
    "model_11_5_13",
    "model_11_5_23",
    "model_11_5_12",
    "model_11_5_10",
    "model_11_5_9",
    "model_11_5_8",
    "model_11_5_7",
    "model_11_5_6",
    "model_11_5_5",
    "model_11_5_4",
    "model_11_5_3",
    "model_11_5_2",
    "model_11_5_1",
    "model_11_5_11",
    "model_11_5_24"
)

# The metric values (columns B..I) that every data row now shares
$metrics = @(
    0.3494677884409869,
    0.09518715775284903,
    -1.151154049613713,
    0.006537356941730965,
    0.7199474573135376,
    1.041117906570435,
    0.3009430766105652,
    0.6928002834320068
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($j = 0; $j -lt $metrics.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $metrics[$j]
    }
}
